$d = $word.ActiveDocument

# 1. Insert the new sentence about categorizing year right after
#    "...same decade. " and before "One of the main concerns"
$d.Content.Find.Execute(
    "same decade. One of the main concerns",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "same decade. We did look at categorizing year as described once we found that treating is as a continuous predictor was not working. One of the main concerns",
    2)

# 2. Reword the PCA interpretability clause
$d.Content.Find.Execute(
    "PCA" + [char]0x2014 + "that it becomes difficult to interpret results-",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "PCA" + [char]0x2014 + "difficulty in interpreting results-",
    2)

# 3. Append "anyway" to the final sentence of the paragraph
$d.Content.Find.Execute(
    "interpretation difficult.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "interpretation difficult anyway.",
    2)
